$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the bold/bordered style
# used by the other header cells (B1:H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-48: pairs of (I0, IF) values for each game row.
$values = @(
    @(7,7),
    @(7,7),
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(7,7),
    @(9,9),
    @(10,10),
    @(7,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(6,6),
    @(8,8),
    @(8,8),
    @(7,7),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(11,11),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(6,6),
    @(4,4),
    @(4,4)
)

$r = 2
foreach ($pair in $values) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r = $r + 1
}
